$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) cells whose new values would otherwise be
# auto-converted to numbers by Excel to remain plain text,
# matching the source data which stores these as text labels.
$textPriceRows = @(5,6,8,9,10,11,12,14,15,16,18,20,21,22,23,25,26,27,28,29,30,31,32,33,34,35,36,37,39,40,41,42,43,45,48,49,50,51)
foreach ($r in $textPriceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Coin / Link swap for rows 48-49 (Aave <-> Mantle reordered)
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"

# Updated Price (D) values
$ws.Range("D2").Value = "29.628.72"
$ws.Range("D3").Value = "1.862.99"
$ws.Range("D5").Value = "245.78"
$ws.Range("D6").Value = "0.7004"
$ws.Range("D8").Value = "0.07756"
$ws.Range("D9").Value = "0.3074"
$ws.Range("D10").Value = "23.69"
$ws.Range("D11").Value = "0.07812"
$ws.Range("D12").Value = "5.169"
$ws.Range("D13").Value = "1.854.20"
$ws.Range("D14").Value = "92.77"
$ws.Range("D15").Value = "0.6939"
$ws.Range("D16").Value = "6.616"
$ws.Range("D17").Value = "29.614.86"
$ws.Range("D18").Value = "0.000008371"
$ws.Range("D19").Value = "2.109.21"
$ws.Range("D20").Value = "243.01"
$ws.Range("D21").Value = "12.80"
$ws.Range("D22").Value = "1.000"
$ws.Range("D23").Value = "7.637"
$ws.Range("D25").Value = "0.1516"
$ws.Range("D26").Value = "8.945"
$ws.Range("D27").Value = "160.00"
$ws.Range("D28").Value = "18.37"
$ws.Range("D29").Value = "1.545"
$ws.Range("D30").Value = "4.272"
$ws.Range("D31").Value = "4.202"
$ws.Range("D32").Value = "1.198"
$ws.Range("D33").Value = "0.05107"
$ws.Range("D34").Value = "0.7871"
$ws.Range("D35").Value = "1.909"
$ws.Range("D36").Value = "1.158"
$ws.Range("D37").Value = "2.693"
$ws.Range("D38").Value = "1.338.49"
$ws.Range("D39").Value = "0.01882"
$ws.Range("D40").Value = "2.740"
$ws.Range("D41").Value = "0.9602"
$ws.Range("D42").Value = "6.022"
$ws.Range("D43").Value = "106.65"
$ws.Range("D45").Value = "9.796"
$ws.Range("D47").Value = "2.011.46"
$ws.Range("D48").Value = "0.5216"
$ws.Range("D49").Value = "65.11"
$ws.Range("D50").Value = "1.790"
$ws.Range("D51").Value = "7.022"

# Updated Volume(1h) (E) values
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +2.68%  "
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +1.69%  "
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("E12").Value = "  +2.56%  "
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("E14").Value = "  +2.94%  "
$ws.Range("E15").Value = "  +3.24%  "
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("E17").Value = "  +2.69%  "
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  +3.16%  "
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("E26").Value = "  +2.80%  "
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("E33").Value = "  +0.23%  "
$ws.Range("E34").Value = "  +3.89%  "
$ws.Range("E35").Value = "  +5.41%  "
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("E38").Value = "  +10.28%  "
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("E40").Value = "  +2.49%  "
$ws.Range("E41").Value = "  +5.04%  "
$ws.Range("E42").Value = "  +14.33%  "
$ws.Range("E43").Value = "  -1.58%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("E45").Value = "  +3.93%  "
$ws.Range("E46").Value = "  +2.85%  "
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("E49").Value = "  +4.16%  "
$ws.Range("E50").Value = "  +3.73%  "
$ws.Range("E51").Value = "  +2.32%  "
